$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 4631930  # row 80: was 5002532.5
$ws.Range("I80").Value = 7814616  # row 80: was 7355008.5
$ws.Range("J80").Value = 2568.818  # row 80: was 3519.875
$ws.Range("K80").Value = 23443848  # row 80: was 22065025.5
$ws.Range("L80").Value = 7706.454000000001  # row 80: was 10559.625
$ws.Range("M80").Value = -23442850  # row 80: was -22064027.5
$ws.Range("N80").Value = -9702.454000000002  # row 80: was -12555.625

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H83").Value = 4631930  # row 83: was 5002532.5
$ws.Range("I83").Value = 7814616  # row 83: was 7355008.5
$ws.Range("J83").Value = 2568.818  # row 83: was 3519.875
$ws.Range("K83").Value = 70331544  # row 83: was 66195076.5
$ws.Range("L83").Value = 23119.362  # row 83: was 31678.875
$ws.Range("M83").Value = -70326552  # row 83: was -66190084.5
$ws.Range("N83").Value = -33103.362  # row 83: was -41662.875

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 1124.5883  # row 112: was 1151.3636
$ws.Range("J112").Value = 1143.129  # row 112: was 1173.2
$ws.Range("L112").Value = 3429.387  # row 112: was 3519.6
$ws.Range("N112").Value = -5645.387  # row 112: was -5735.6

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 39251.953  # row 132: was 40546.13
$ws.Range("I132").Value = 44118.383  # row 132: was 46635.133
$ws.Range("J132").Value = 9512.666999999999  # row 132: was 8883.299999999999
$ws.Range("K132").Value = 132355.149  # row 132: was 139905.399
$ws.Range("L132").Value = 28538.001  # row 132: was 26649.9
$ws.Range("M132").Value = -129825.149  # row 132: was -137375.399
$ws.Range("N132").Value = -33598.001  # row 132: was -31709.9

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2262.4424  # row 138: was 2283.2942
$ws.Range("J138").Value = 2635.6216  # row 138: was 2675.5278
$ws.Range("L138").Value = 7906.864799999999  # row 138: was 8026.5834
$ws.Range("N138").Value = -18186.8648  # row 138: was -18306.5834

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H28").Value = 1905  # row 28: was 3900
$ws.Range("I28").Value = 1905  # row 28: was 3000
$ws.Range("J28").Value = 0  # row 28: was 4800
$ws.Range("K28").Value = 1905  # row 28: was 3000
$ws.Range("L28").Value = 0  # row 28: was 4800
$ws.Range("M28").Value = -1713  # row 28: was -2808
$ws.Range("N28").ClearContents()  # row 28: was -5184

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 4465899.5  # row 74: was 4033737.5
$ws.Range("I74").Value = 5435790.5  # row 74: was 4630550.5
$ws.Range("J74").Value = 4400  # row 74: was 5250
$ws.Range("K74").Value = 5435790.5  # row 74: was 4630550.5
$ws.Range("L74").Value = 4400  # row 74: was 5250
$ws.Range("M74").Value = -5434916.5  # row 74: was -4629676.5
$ws.Range("N74").Value = -6148  # row 74: was -6998

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 4465899.5  # row 77: was 4033737.5
$ws.Range("I77").Value = 5435790.5  # row 77: was 4630550.5
$ws.Range("J77").Value = 4400  # row 77: was 5250
$ws.Range("K77").Value = 27178952.5  # row 77: was 23152752.5
$ws.Range("L77").Value = 22000  # row 77: was 26250
$ws.Range("M77").Value = -27174584.5  # row 77: was -23148384.5
$ws.Range("N77").Value = -30736  # row 77: was -34986

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 1181.8334  # row 97: was 1207.8235
$ws.Range("I97").Value = 1181.8334  # row 97: was 1207.8235
$ws.Range("K97").Value = 1181.8334  # row 97: was 1207.8235
$ws.Range("M97").Value = -685.8334  # row 97: was -711.8235

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H99").Value = 1905  # row 99: was 3900
$ws.Range("I99").Value = 1905  # row 99: was 3000
$ws.Range("J99").Value = 0  # row 99: was 4800
$ws.Range("K99").Value = 1905  # row 99: was 3000
$ws.Range("L99").Value = 0  # row 99: was 4800
$ws.Range("M99").Value = 1090  # row 99: was -5
$ws.Range("N99").ClearContents()  # row 99: was -10790

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 351545.16  # row 132: was 378129.47
$ws.Range("I132").Value = 486227.53  # row 132: was 543351.3
$ws.Range("K132").Value = 1458682.59  # row 132: was 1630053.9
$ws.Range("M132").Value = -1456152.59  # row 132: was -1627523.9

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1662.8  # row 86: was 1701.1724
$ws.Range("I86").Value = 1515.1765  # row 86: was 1575.5
$ws.Range("K86").Value = 1515.1765  # row 86: was 1575.5
$ws.Range("M86").Value = -392.1765  # row 86: was -452.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 1662.8  # row 89: was 1701.1724
$ws.Range("I89").Value = 1515.1765  # row 89: was 1575.5
$ws.Range("K89").Value = 7575.8825  # row 89: was 7877.5
$ws.Range("M89").Value = -1959.8825  # row 89: was -2261.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1944  # row 105: was 0
$ws.Range("I105").Value = 1944  # row 105: was 0
$ws.Range("K105").Value = 1944  # row 105: was 0
$ws.Range("M105").Value = -197  # row 105: was None

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1919968.8  # row 134: was 2018256.4
$ws.Range("J134").Value = 767575  # row 134: was 8331.666999999999
$ws.Range("L134").Value = 2302725  # row 134: was 24995.001
$ws.Range("N134").Value = -2307795  # row 134: was -30065.001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 401830.06  # row 58: was 415201.12
$ws.Range("I58").Value = 589302.6  # row 58: was 618732.9
$ws.Range("K58").Value = 589302.6  # row 58: was 618732.9
$ws.Range("M58").Value = -589099.6  # row 58: was -618529.9

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H88").Value = 46419.57  # row 88: was 49346.855
$ws.Range("J88").Value = 46989.5  # row 88: was 50404.668
$ws.Range("L88").Value = 46989.5  # row 88: was 50404.668
$ws.Range("N88").Value = -47801.5  # row 88: was -51216.668

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H91").Value = 46419.57  # row 91: was 49346.855
$ws.Range("J91").Value = 46989.5  # row 91: was 50404.668
$ws.Range("L91").Value = 46989.5  # row 91: was 50404.668
$ws.Range("N91").Value = -49797.5  # row 91: was -53212.668

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 75016600  # row 132: was 75016680
$ws.Range("I132").Value = 111134264  # row 132: was 111134400
$ws.Range("K132").Value = 333402792  # row 132: was 333403200
$ws.Range("M132").Value = -333400262  # row 132: was -333400670

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 22004.176  # row 134: was 23357.875
$ws.Range("I134").Value = 27352.154  # row 134: was 29602.75
$ws.Range("K134").Value = 82056.462  # row 134: was 88808.25
$ws.Range("M134").Value = -79521.462  # row 134: was -86273.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 401830.06  # row 136: was 415201.12
$ws.Range("I136").Value = 589302.6  # row 136: was 618732.9
$ws.Range("K136").Value = 1767907.8  # row 136: was 1856198.7
$ws.Range("M136").Value = -1765357.8  # row 136: was -1853648.7

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 458.75  # row 8: was 431.6
$ws.Range("I8").Value = 458.75  # row 8: was 431.6
$ws.Range("K8").Value = 1376.25  # row 8: was 1294.8
$ws.Range("M8").Value = -1237.25  # row 8: was -1155.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H82").Value = 2800  # row 82: was 5188
$ws.Range("I82").Value = 2800  # row 82: was 2797.5
$ws.Range("J82").Value = 0  # row 82: was 9969
$ws.Range("K82").Value = 8400  # row 82: was 8392.5
$ws.Range("L82").Value = 0  # row 82: was 29907
$ws.Range("M82").Value = -7994  # row 82: was -7986.5
$ws.Range("N82").ClearContents()  # row 82: was -30719

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H85").Value = 2800  # row 85: was 5188
$ws.Range("I85").Value = 2800  # row 85: was 2797.5
$ws.Range("J85").Value = 0  # row 85: was 9969
$ws.Range("K85").Value = 8400  # row 85: was 8392.5
$ws.Range("L85").Value = 0  # row 85: was 29907
$ws.Range("M85").Value = -6996  # row 85: was -6988.5
$ws.Range("N85").ClearContents()  # row 85: was -32715

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 243  # row 107: was 247
$ws.Range("I107").Value = 242  # row 107: was 247
$ws.Range("J107").Value = 244  # row 107: was 0
$ws.Range("K107").Value = 726  # row 107: was 741
$ws.Range("L107").Value = 732  # row 107: was 0
$ws.Range("M107").Value = 1194  # row 107: was 1179
$ws.Range("N107").Value = -4572  # row 107: was None

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H108").Value = 14504  # row 108: was 7316.7
$ws.Range("I108").Value = 400  # row 108: was 174.5
$ws.Range("K108").Value = 1200  # row 108: was 523.5
$ws.Range("M108").Value = 1680  # row 108: was 2356.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H118").Value = 909.6667  # row 118: was 1114.5
$ws.Range("I118").Value = 764.5  # row 118: was 1029
$ws.Range("K118").Value = 2293.5  # row 118: was 3087
$ws.Range("M118").Value = -1050.5  # row 118: was -1844

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 3081.8333  # row 137: was 3310.5
$ws.Range("J137").Value = 3756.5625  # row 137: was 4212.2856
$ws.Range("L137").Value = 11269.6875  # row 137: was 12636.8568
$ws.Range("N137").Value = -21469.6875  # row 137: was -22836.8568

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 45906.5  # row 107: was 52421.855
$ws.Range("J107").Value = 2926.25  # row 107: was 3802
$ws.Range("L107").Value = 2926.25  # row 107: was 3802
$ws.Range("N107").Value = -6766.25  # row 107: was -7642

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 5420.6875  # row 122: was 5060.5
$ws.Range("I122").Value = 3196.4546  # row 122: was 3219.5715
$ws.Range("J122").Value = 10314  # row 122: was 11503.75
$ws.Range("K122").Value = 9589.363799999999  # row 122: was 9658.7145
$ws.Range("L122").Value = 30942  # row 122: was 34511.25
$ws.Range("M122").Value = -7139.363799999999  # row 122: was -7208.7145
$ws.Range("N122").Value = -35842  # row 122: was -39411.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 14677735  # row 132: was 15824284
$ws.Range("I132").Value = 22500976  # row 132: was 24696070
$ws.Range("J132").Value = 9158.375  # row 132: was 9360.305
$ws.Range("K132").Value = 67502928  # row 132: was 74088210
$ws.Range("L132").Value = 27475.125  # row 132: was 28080.915
$ws.Range("M132").Value = -67500398  # row 132: was -74085680
$ws.Range("N132").Value = -32535.125  # row 132: was -33140.915

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3522.1667  # row 7: was 4599.885
$ws.Range("I7").Value = 3462.4375  # row 7: was 3224.95
$ws.Range("J7").Value = 4000  # row 7: was 9183
$ws.Range("K7").Value = 3462.4375  # row 7: was 3224.95
$ws.Range("L7").Value = 4000  # row 7: was 9183
$ws.Range("M7").Value = -3350.4375  # row 7: was -3112.95
$ws.Range("N7").Value = -4224  # row 7: was -9407

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 3522.1667  # row 126: was 4599.885
$ws.Range("I126").Value = 3462.4375  # row 126: was 3224.95
$ws.Range("J126").Value = 4000  # row 126: was 9183
$ws.Range("K126").Value = 10387.3125  # row 126: was 9674.849999999999
$ws.Range("L126").Value = 12000  # row 126: was 27549
$ws.Range("M126").Value = -7917.3125  # row 126: was -7204.849999999999
$ws.Range("N126").Value = -16940  # row 126: was -32489

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H49").Value = 0  # row 49: was 42000
$ws.Range("I49").Value = 0  # row 49: was 42000
$ws.Range("K49").Value = 0  # row 49: was 42000
$ws.Range("M49").ClearContents()  # row 49: was -41770

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1897.0741  # row 107: was 2020.76
$ws.Range("I107").Value = 1121.7727  # row 107: was 1198.85
$ws.Range("K107").Value = 3365.3181  # row 107: was 3596.55
$ws.Range("M107").Value = -1445.3181  # row 107: was -1676.55

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1581.5358  # row 113: was 1530.1034
$ws.Range("I113").Value = 442.92856  # row 113: was 419.4
$ws.Range("K113").Value = 1328.78568  # row 113: was 1258.2
$ws.Range("M113").Value = 841.21432  # row 113: was 911.8000000000002

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 3557.3103  # row 122: was 3081.6
$ws.Range("I122").Value = 3416.6785  # row 122: was 2951.7942
$ws.Range("K122").Value = 10250.0355  # row 122: was 8855.382599999999
$ws.Range("M122").Value = -7800.0355  # row 122: was -6405.382599999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H125").Value = 80715  # row 125: was 30000
$ws.Range("J125").Value = 80715  # row 125: was 30000
$ws.Range("L125").Value = 80715  # row 125: was 30000
$ws.Range("N125").Value = -90555  # row 125: was -39840

Write-Host "done applying 38 row updates"
